# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- 1) Reorder country names: Hungria now comes before Senegal ---
# Row 87 was "Senegal" with its stats; Row 88 was "Hungria" with its stats.
# After the update, Row 87 becomes "Hungria" (with NEW/updated stats) and
# Row 88 becomes "Senegal" (keeping its previous stats).
$ws.Range("A87").Value = "Hungria"
$ws.Range("A88").Value = "Senegal"

# --- 2) Update numeric data ---
# Columns: A=Pais B=Casos totales C=Nuevos casos D=Casos activos E=Recuperados F=Casos criticos G=Muertes hoy H=Muertes

# Row 7 - Rusia
$ws.Range("B7").Value = 1085281
$ws.Range("C7").Value = 5762
$ws.Range("D7").Value = 895868
$ws.Range("E7").Value = 170352
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 144
$ws.Range("H7").Value = 19061

# Row 62 - Armenia
$ws.Range("B62").Value = 46671
$ws.Range("C62").Value = 295
$ws.Range("D62").Value = 42231
$ws.Range("E62").Value = 3515
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 925

# Row 87 - Hungria (new data)
$ws.Range("B87").Value = 15170
$ws.Range("C87").Value = 710
$ws.Range("D87").Value = 4227
$ws.Range("E87").Value = 10280
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 9
$ws.Range("H87").Value = 663

# Row 88 - Senegal (previous data, now shifted down one rank)
$ws.Range("B88").Value = 14568
$ws.Range("C88").Value = 0
$ws.Range("D88").Value = 10756
$ws.Range("E88").Value = 3513
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 299

# Row 143 - Georgia
$ws.Range("B143").Value = 2937
$ws.Range("C143").Value = 179
$ws.Range("D143").Value = 1422
$ws.Range("E143").Value = 1496
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 19

# Row 160 - Letonia
$ws.Range("B160").Value = 1494
$ws.Range("C160").Value = 8
$ws.Range("D160").Value = 1248
$ws.Range("E160").Value = 210
$ws.Range("F160").Value = 0
$ws.Range("G160").Value = 1
$ws.Range("H160").Value = 36

# Row 176 - Taiwan
$ws.Range("B176").Value = 503
$ws.Range("C176").Value = 3
$ws.Range("D176").Value = 478
$ws.Range("E176").Value = 18
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 7
